{"js": "// The three \"git ...\" paragraphs were originally split into two runs each\n// (e.g. \"git\" / \" add <file name>\") bracketed by spellcheck/grammar-check\n// <w:proofErr/> markers. The edit removes those markers and merges each\n// paragraph's text back into a single run, then appends one more\n// instruction paragraph after \"git push\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Helper: rebuild a paragraph's content as a single run (dropping any\n// <w:proofErr/> marks and extra run splits) while preserving its <w:pPr/>\n// and original rsid attributes.\nfunction replaceParagraphText(paragraph, text) {\n  const escaped = text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n  const ooxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    '<w:p w:rsidR=\"0032294A\" w:rsidRDefault=\"0032294A\" w:rsidP=\"0032294A\">' +\n    '<w:pPr><w:spacing w:after=\"0\"/></w:pPr>' +\n    \"<w:r>\" +\n    \"<w:t>\" +\n    escaped +\n    \"</w:t>\" +\n    \"</w:r>\" +\n    \"</w:p>\" +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\";\n  paragraph.getRange().insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\nconst gitAddPara = paragraphs.items[2];\nconst gitCommitPara = paragraphs.items[3];\nconst gitPushPara = paragraphs.items[4];\n\nreplaceParagraphText(gitAddPara, \"git add <file name>\");\nreplaceParagraphText(gitCommitPara, \"git commit \\u2013m \\u201c<write what is updated>\\u201d\");\nreplaceParagraphText(gitPushPara, \"git push\");\nawait context.sync();\n\n// Add the new trailing instruction paragraph.\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastGitPara = paragraphs.items[4];\nlastGitPara.insertParagraph(\n  \"(user name & password will be asked)\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# The three \"git ...\" paragraphs were originally split into two runs each\n# (e.g. \"git\" / \" add <file name>\") bracketed by spellcheck/grammar-check\n# <w:proofErr/> markers. This edit removes those markers and merges each\n# paragraph's text back into a single run, then appends one more\n# instruction paragraph after \"git push\".\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText($paragraph, [string]$text) {\n    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n    $ooxml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p w:rsidR=\"0032294A\" w:rsidRDefault=\"0032294A\" w:rsidP=\"0032294A\">' +\n        '<w:pPr><w:spacing w:after=\"0\"/></w:pPr>' +\n        '<w:r><w:t>' + $escaped + '</w:t></w:r>' +\n        '</w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n    # Range.InsertXML replaces this paragraph's contents wholesale, which\n    # drops the stray <w:proofErr/> markers and collapses the two runs\n    # (\"git\" + \" add <file name>\") into the single run seen above.\n    $paragraph.Range.InsertXML($ooxml)\n}\n\n$gitAddPara = $d.Paragraphs(3)\nSet-ParagraphText $gitAddPara \"git add <file name>\"\n\n$gitCommitPara = $d.Paragraphs(4)\nSet-ParagraphText $gitCommitPara \"git commit \u2013m \u201c<write what is updated>\u201d\"\n\n# Insert the new trailing paragraph BEFORE rewriting \"git push\" so that\n# paragraph is no longer the very last one in the body when InsertXML\n# runs on it (doing it after would leave a spurious empty paragraph).\n$gitPushPara = $d.Paragraphs(5)\n$gitPushPara.Range.InsertParagraphAfter()\n\n$gitPushPara = $d.Paragraphs(5)\nSet-ParagraphText $gitPushPara \"git push\"\n\n$newPara = $d.Paragraphs(6)\n$newPara.Range.Text = \"(user name & password will be asked)\"\n"}
